# Apply the "Added on 20, 2017" edit:
#  - Sheet1: update the counter cell A8, move the remembered selection to A9,
#            and stop being the active tab.
#  - Sheet2: becomes the active tab; gains two new label/value rows and is
#            selected at E2; the label column is auto-widened for the new text.
#  - workbook-level active tab flips from Sheet1 (0) to Sheet2 (1).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: bump the running counter in A8 and leave the selection on A9 ---
$ws1.Range("A8").Value = 3013685163
$ws1.Range("A9").Select() | Out-Null

# --- Sheet2: new rows of data ---
$ws2.Range("A1").Value = "Schedule Number Counter"
$ws2.Range("B1").Value = 132
$ws2.Range("A2").Value = "Data Recorder Index"
$ws2.Range("B2").Value = 12

# Size column A to fit the new label text.
$ws2.Columns.Item(1).ColumnWidth = 24

# Sheet2 becomes the active sheet/tab, selected at E2.
$ws2.Activate() | Out-Null
$ws2.Range("E2").Select() | Out-Null
